$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 38
$ws.Range("H38").Value = 2622.3
$ws.Range("J38").Value = 6485
$ws.Range("L38").Value = 19455
$ws.Range("N38").Value = -20199

# Row 51
$ws.Range("H51").Value = 2937.5
$ws.Range("I51").Value = 2750
$ws.Range("J51").Value = 3000
$ws.Range("K51").Value = 2750
$ws.Range("L51").Value = 3000
$ws.Range("M51").Value = -2266
$ws.Range("N51").Value = -3968

# Row 97
$ws.Range("H97").Value = 565.6
$ws.Range("J97").Value = 565.6
$ws.Range("L97").Value = 1696.8
$ws.Range("N97").Value = -2688.8

# Row 98
$ws.Range("H98").Value = 568.1111
$ws.Range("I98").Value = 489.25
$ws.Range("K98").Value = 489.25
$ws.Range("M98").Value = 1008.75

# Row 122
$ws.Range("H122").Value = 568.1111
$ws.Range("I122").Value = 489.25
$ws.Range("K122").Value = 1467.75
$ws.Range("M122").Value = 982.25

# Row 132
$ws.Range("H132").Value = 8223
$ws.Range("I132").Value = 8167.5
$ws.Range("J132").Value = 9000
$ws.Range("K132").Value = 24502.5
$ws.Range("L132").Value = 27000
$ws.Range("M132").Value = -21972.5
$ws.Range("N132").Value = -32060

# Row 135
$ws.Range("H135").Value = 659.4545000000001
$ws.Range("I135").Value = 659.4545000000001
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 5935.0905
$ws.Range("L135").Value = 0
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -3400.0905

# Row 137
$ws.Range("H137").Value = 2129.3572
$ws.Range("I137").Value = 866.1429000000001
$ws.Range("J137").Value = 3392.5715
$ws.Range("K137").Value = 2598.4287
$ws.Range("L137").Value = 10177.7145
$ws.Range("M137").Value = -48.42870000000039
$ws.Range("N137").Value = -15277.7145

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 1795.1
$ws.Range("I61").Value = 1270.6666
$ws.Range("J61").Value = 2581.75
$ws.Range("K61").Value = 1270.6666
$ws.Range("L61").Value = 2581.75
$ws.Range("M61").Value = -1058.6666
$ws.Range("N61").Value = -3005.75

# Row 136
$ws.Range("H136").Value = 1795.1
$ws.Range("I136").Value = 1270.6666
$ws.Range("J136").Value = 2581.75
$ws.Range("K136").Value = 3811.9998
$ws.Range("L136").Value = 7745.25
$ws.Range("M136").Value = -1261.9998
$ws.Range("N136").Value = -12845.25

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 4009.5
$ws.Range("I99").Value = 4009.5
$ws.Range("K99").Value = 4009.5
$ws.Range("M99").Value = -2511.5

# Row 106
$ws.Range("H106").Value = 50166.668
$ws.Range("J106").Value = 50166.668
$ws.Range("L106").Value = 50166.668
$ws.Range("N106").Value = -52690.668

# Row 107
$ws.Range("H107").Value = 5300.6665
$ws.Range("I107").Value = 1541.2
$ws.Range("J107").Value = 10000
$ws.Range("K107").Value = 1541.2
$ws.Range("L107").Value = 10000
$ws.Range("M107").Value = 378.8
$ws.Range("N107").Value = -13840

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 952.2
$ws.Range("I16").Value = 502.75
$ws.Range("K16").Value = 502.75
$ws.Range("M16").Value = -215.75

# Row 31
$ws.Range("H31").Value = 6029.2
$ws.Range("I31").Value = 1984.5
$ws.Range("J31").Value = 8725.666999999999
$ws.Range("K31").Value = 1984.5
$ws.Range("L31").Value = 8725.666999999999
$ws.Range("M31").Value = -1689.5
$ws.Range("N31").Value = -9315.666999999999

# Row 34
$ws.Range("H34").Value = 6029.2
$ws.Range("I34").Value = 1984.5
$ws.Range("J34").Value = 8725.666999999999
$ws.Range("K34").Value = 1984.5
$ws.Range("L34").Value = 8725.666999999999
$ws.Range("M34").Value = -1782.5
$ws.Range("N34").Value = -9129.666999999999

# Row 58
$ws.Range("H58").Value = 2326.6316
$ws.Range("I58").Value = 1093.7142
$ws.Range("K58").Value = 1093.7142
$ws.Range("M58").Value = -890.7141999999999

# Row 68
$ws.Range("H68").Value = 66320
$ws.Range("J68").Value = 66320
$ws.Range("L68").Value = 66320
$ws.Range("N68").Value = -67818

# Row 71
$ws.Range("H71").Value = 66320
$ws.Range("J71").Value = 66320
$ws.Range("L71").Value = 198960
$ws.Range("N71").Value = -206448

# Row 107
$ws.Range("H107").Value = 2032
$ws.Range("I107").Value = 550
$ws.Range("J107").Value = 4996
$ws.Range("K107").Value = 550
$ws.Range("L107").Value = 4996
$ws.Range("M107").Value = 1370
$ws.Range("N107").Value = -8836

# Row 113
$ws.Range("H113").Value = 952.2
$ws.Range("I113").Value = 502.75
$ws.Range("K113").Value = 502.75
$ws.Range("M113").Value = 1667.25

# Row 117
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").ClearContents()
$ws.Range("N117").Value = 0

# Row 132
$ws.Range("H132").Value = 1949.8462
$ws.Range("I132").Value = 1949.8462
$ws.Range("K132").Value = 5849.5386
$ws.Range("M132").Value = -3319.5386

# Row 134
$ws.Range("H134").Value = 2221.8333
$ws.Range("I134").Value = 2221.8333
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 6665.499899999999
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -4130.499899999999

# Row 136
$ws.Range("H136").Value = 2326.6316
$ws.Range("I136").Value = 1093.7142
$ws.Range("K136").Value = 3281.1426
$ws.Range("M136").Value = -731.1425999999997

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1088.75
$ws.Range("I5").Value = 908.3333
$ws.Range("K5").Value = 2724.9999
$ws.Range("M5").Value = -2612.9999

# Row 38
$ws.Range("H38").Value = 452.20834
$ws.Range("I38").Value = 431.78946
$ws.Range("K38").Value = 1295.36838
$ws.Range("M38").Value = -948.3683800000001

# Row 43
$ws.Range("H43").Value = 1000
$ws.Range("I43").Value = 1000
$ws.Range("K43").Value = 3000
$ws.Range("M43").Value = -2886

# Row 55
$ws.Range("H55").Value = 7833.143
$ws.Range("J55").Value = 7833.143
$ws.Range("L55").Value = 23499.429
$ws.Range("N55").Value = -23853.429

# Row 87
$ws.Range("H87").Value = 1500
$ws.Range("I87").Value = 1500
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 4500
$ws.Range("L87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("N87").Value = -3252

# Row 90
$ws.Range("H90").Value = 1500
$ws.Range("I90").Value = 1500
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 13500
$ws.Range("L90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("N90").Value = -7260

# Row 135
$ws.Range("H135").Value = 1088.75
$ws.Range("I135").Value = 908.3333
$ws.Range("K135").Value = 8174.9997
$ws.Range("M135").Value = -5639.9997

# Row 137
$ws.Range("H137").Value = 3321.5
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 3321.5
$ws.Range("K137").Value = 0
$ws.Range("L137").ClearContents()
$ws.Range("M137").Value = 9964.5
$ws.Range("N137").Value = -20164.5

$ws = $wb.Worksheets.Item("GSM")
# Row 3
$ws.Range("H3").Value = 25648898
$ws.Range("J3").Value = 40001332
$ws.Range("L3").Value = 40001332
$ws.Range("N3").Value = -40001564

# Row 113
$ws.Range("H113").Value = 4451.6665
$ws.Range("I113").Value = 1632.3334
$ws.Range("J113").Value = 7271
$ws.Range("K113").Value = 1632.3334
$ws.Range("L113").Value = 7271
$ws.Range("M113").Value = 537.6666
$ws.Range("N113").Value = -11611

# Row 126
$ws.Range("H126").Value = 2601.4614
$ws.Range("I126").Value = 1993.25
$ws.Range("K126").Value = 5979.75
$ws.Range("M126").Value = -3509.75

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 8219.714
$ws.Range("I7").Value = 7787.8
$ws.Range("J7").Value = 9299.5
$ws.Range("K7").Value = 7787.8
$ws.Range("L7").Value = 9299.5
$ws.Range("M7").Value = -7675.8
$ws.Range("N7").Value = -9523.5

# Row 46
$ws.Range("H46").Value = 6166.4165
$ws.Range("I46").Value = 1749.25
$ws.Range("J46").Value = 8375
$ws.Range("K46").Value = 1749.25
$ws.Range("L46").Value = 8375
$ws.Range("M46").Value = -1561.25
$ws.Range("N46").Value = -8751

# Row 82
$ws.Range("H82").Value = 5964.143
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()

# Row 85
$ws.Range("H85").Value = 5964.143
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()

# Row 126
$ws.Range("H126").Value = 8219.714
$ws.Range("I126").Value = 7787.8
$ws.Range("J126").Value = 9299.5
$ws.Range("K126").Value = 23363.4
$ws.Range("L126").Value = 27898.5
$ws.Range("M126").Value = -20893.4
$ws.Range("N126").Value = -32838.5

# Row 132
$ws.Range("H132").Value = 1647.3334
$ws.Range("I132").Value = 996.8
$ws.Range("J132").Value = 4900
$ws.Range("K132").Value = 2990.4
$ws.Range("L132").Value = 14700
$ws.Range("M132").Value = -460.3999999999996
$ws.Range("N132").Value = -19760

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 1084.9375
$ws.Range("I132").Value = 1084.9375
$ws.Range("K132").Value = 3254.8125
$ws.Range("M132").Value = -724.8125
